$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.629.94"
$ws.Range("E2").Value = "  -0.61%  "

$ws.Range("D3").Value = "3.444.14"
$ws.Range("E3").Value = "  -2.42%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "590.73"
$ws.Range("E5").Value = "  -1.69%  "

$ws.Range("D6").Value = "178.79"
$ws.Range("E6").Value = "  -2.68%  "

$ws.Range("E7").Value = "  +1.72%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "3.441.41"
$ws.Range("E9").Value = "  -2.47%  "

$ws.Range("E10").Value = "  -1.31%  "

$ws.Range("D11").Value = "6.97"
$ws.Range("E11").Value = "  -2.54%  "

$ws.Range("D12").Value = "0.427"
$ws.Range("E12").Value = "  -3.57%  "

$ws.Range("D13").Value = "4.039.09"
$ws.Range("E13").Value = "  -2.48%  "

$ws.Range("D14").Value = "32.02"
$ws.Range("E14").Value = "  -1.53%  "

$ws.Range("E15").Value = "  -1.22%  "

$ws.Range("D16").Value = "67.601.34"
$ws.Range("E16").Value = "  -0.62%  "

$ws.Range("D17").Value = "0.0000176"
$ws.Range("E17").Value = "  -3.65%  "

$ws.Range("D18").Value = "3.444.04"
$ws.Range("E18").Value = "  -2.58%  "

$ws.Range("D19").Value = "6.13"
$ws.Range("E19").Value = "  -4.55%  "

$ws.Range("D20").Value = "13.98"
$ws.Range("E20").Value = "  -6.73%  "

$ws.Range("D21").Value = "389.16"
$ws.Range("E21").Value = "  -2.57%  "

$ws.Range("D22").Value = "7.84"
$ws.Range("E22").Value = "  -3.74%  "

$ws.Range("E23").Value = "  +1.91%  "

$ws.Range("E24").Value = "  -0.29%  "

$ws.Range("D25").Value = "0.532"
$ws.Range("E25").Value = "  -2.80%  "

$ws.Range("D26").Value = "71.34"
$ws.Range("E26").Value = "  -3.10%  "

$ws.Range("E27").Value = "  -5.53%  "

$ws.Range("D28").Value = "10.24"
$ws.Range("E28").Value = "  -4.98%  "

$ws.Range("E29").Value = "  -2.27%  "

$ws.Range("E30").Value = "  +0.33%  "

$ws.Range("D31").Value = "6.04"
$ws.Range("E31").Value = "  -4.32%  "

$ws.Range("E32").Value = "  -2.20%  "

$ws.Range("E33").Value = "  -6.29%  "

$ws.Range("D34").Value = "23.22"
$ws.Range("E34").Value = "  -3.81%  "

$ws.Range("E35").Value = "  -0.12%  "

$ws.Range("D36").Value = "7.21"
$ws.Range("E36").Value = "  -3.92%  "

$ws.Range("E37").Value = "  -7.99%  "

$ws.Range("D38").Value = "161.06"
$ws.Range("E38").Value = "  -1.77%  "

$ws.Range("D39").Value = "0.883"
$ws.Range("E39").Value = "  +0.20%  "

$ws.Range("E40").Value = "  -6.09%  "

$ws.Range("D41").Value = "2.73"
$ws.Range("E41").Value = "  -2.01%  "

$ws.Range("E42").Value = "  -3.87%  "

$ws.Range("D43").Value = "6.59"
$ws.Range("E43").Value = "  -8.40%  "

$ws.Range("D44").Value = "25.76"
$ws.Range("E44").Value = "  -5.16%  "

$ws.Range("D45").Value = "0.0712"
$ws.Range("E45").Value = "  -4.06%  "

$ws.Range("D46").Value = "25.97"
$ws.Range("E46").Value = "  -6.17%  "

$ws.Range("D47").Value = "2.696.25"
$ws.Range("E47").Value = "  -6.59%  "

$ws.Range("D48").Value = "41.14"
$ws.Range("E48").Value = "  -3.01%  "

$ws.Range("D49").Value = "0.0296"
$ws.Range("E49").Value = "  -3.82%  "

$ws.Range("D50").Value = "324.09"
$ws.Range("E50").Value = "  -8.01%  "

$ws.Range("E51").Value = "  -5.22%  "
